$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute("F&G Life (A Fidelity National Financial company)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $target = $find.Parent
    $target.Collapse(0)
    $target.InsertAfter(" June 2022 through Present")
    $target.Font.Name = "Arial"
    $target.Font.NameBi = "Arial"
    $target.Font.Bold = $true
    $target.Font.Italic = $true
}
